# "Fruta / hortaliza, semanal" - refresh the weekly Espárragos price rows.
# The data rows (2-19) for columns D,H,I,J,K,L,M,N,O,P get re-shuffled onto
# new dates/rows (columns A,B,C,E,F,G,Q,R are identical across every row and
# stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('D','H','I','J','K','L','M','N','O','P')

# target row -> source row (i.e. new row $target gets old row $source's data)
$map = @{
    2 = 5
    3 = 17
    4 = 8
    5 = 9
    6 = 12
    7 = 10
    8 = 11
    9 = 16
    10 = 18
    11 = 7
    12 = 3
    13 = 2
    14 = 4
    15 = 14
    16 = 6
    17 = 19
    18 = 15
    19 = 13
}

# Snapshot every source row's values first (Value2 keeps numbers numeric and
# strings as plain strings) so overwriting earlier rows doesn't clobber data
# that a later row still needs to read.
$orig = @{}
for ($r = 2; $r -le 19; $r++) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

foreach ($target in $map.Keys) {
    $source = $map[$target]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $orig[$source][$c]
    }
}
